$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 -> new values (previously held by row 4)
$ws.Range("D3").Value = 45008
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("S3").Value = 3500

# Row 4 -> new values (previously held by row 5, date unchanged)
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("S4").Value = 3000

# Row 5 -> new values (previously held by row 3)
$ws.Range("D5").Value = 44995
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 5500
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 5750
$ws.Range("S5").Value = 2875
